# Update NATMI TPM-derived LR-pair metrics (Mmp13-Lrp1) to reflect recomputed values.
# Only numeric value cells change; row/column layout and headers are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.6375143333333333
$ws.Range("H2").Value = 1.912543
$ws.Range("I2").Value = 0.005637788761748074
$ws.Range("J2").Value = 0.005637788761748073
$ws.Range("M2").Value = 3.795192333333334
$ws.Range("N2").Value = 11.385577
$ws.Range("O2").Value = 0.01044213755712683
$ws.Range("P2").Value = 0.01044213755712683
$ws.Range("Q2").Value = 2.419489510256778
$ws.Range("R2").Value = 21.775405592311
$ws.Range("S2").Value = 0.00005887056576819715
$ws.Range("T2").Value = 0.00005887056576819713

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.6375143333333333
$ws.Range("H3").Value = 1.912543
$ws.Range("I3").Value = 0.005637788761748074
$ws.Range("J3").Value = 0.005637788761748073
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.6696287328350964
$ws.Range("P3").Value = 0.6696287328350964
$ws.Range("Q3").Value = 155.1559425450476
$ws.Range("R3").Value = 1396.403482905428
$ws.Range("S3").Value = 0.00377522534452131
$ws.Range("T3").Value = 0.003775225344521309

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.6375143333333333
$ws.Range("H4").Value = 1.912543
$ws.Range("I4").Value = 0.005637788761748074
$ws.Range("J4").Value = 0.005637788761748073
$ws.Range("M4").Value = 29.801371
$ws.Range("N4").Value = 89.404113
$ws.Range("O4").Value = 0.08199584844219236
$ws.Range("P4").Value = 0.08199584844219235
$ws.Range("Q4").Value = 18.99880116548433
$ws.Range("R4").Value = 170.989210489359
$ws.Range("S4").Value = 0.0004622752728573904
$ws.Range("T4").Value = 0.0004622752728573903

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.6375143333333333
$ws.Range("H5").Value = 1.912543
$ws.Range("I5").Value = 0.005637788761748074
$ws.Range("J5").Value = 0.005637788761748073
$ws.Range("M5").Value = 86.47679266666667
$ws.Range("N5").Value = 259.430378
$ws.Range("O5").Value = 0.2379332811655844
$ws.Range("P5").Value = 0.2379332811655844
$ws.Range("Q5").Value = 55.13019482569489
$ws.Range("R5").Value = 496.171753431254
$ws.Range("S5").Value = 0.001341417578601176
$ws.Range("T5").Value = 0.001341417578601176

# Row 6
$ws.Range("I6").Value = 0.9321386591533842
$ws.Range("J6").Value = 0.9321386591533841
$ws.Range("M6").Value = 3.795192333333334
$ws.Range("N6").Value = 11.385577
$ws.Range("O6").Value = 0.01044213755712683
$ws.Range("P6").Value = 0.01044213755712683
$ws.Range("Q6").Value = 400.0326729565414
$ws.Range("R6").Value = 3600.294056608872
$ws.Range("S6").Value = 0.009733520101195402
$ws.Range("T6").Value = 0.0097335201011954

# Row 7
$ws.Range("I7").Value = 0.9321386591533842
$ws.Range("J7").Value = 0.9321386591533841
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.6696287328350964
$ws.Range("P7").Value = 0.6696287328350964
$ws.Range("R7").Value = 230878.0490365499
$ws.Range("S7").Value = 0.6241868291554865
$ws.Range("T7").Value = 0.6241868291554864

# Row 8
$ws.Range("I8").Value = 0.9321386591533842
$ws.Range("J8").Value = 0.9321386591533841
$ws.Range("M8").Value = 29.801371
$ws.Range("N8").Value = 89.404113
$ws.Range("O8").Value = 0.08199584844219236
$ws.Range("P8").Value = 0.08199584844219235
$ws.Range("Q8").Value = 3141.216848008552
$ws.Range("R8").Value = 28270.95163207697
$ws.Range("S8").Value = 0.0764315002230493
$ws.Range("T8").Value = 0.07643150022304927

# Row 9
$ws.Range("I9").Value = 0.9321386591533842
$ws.Range("J9").Value = 0.9321386591533841
$ws.Range("M9").Value = 86.47679266666667
$ws.Range("N9").Value = 259.430378
$ws.Range("O9").Value = 0.2379332811655844
$ws.Range("P9").Value = 0.2379332811655844
$ws.Range("Q9").Value = 9115.096016430778
$ws.Range("R9").Value = 82035.86414787702
$ws.Range("S9").Value = 0.221786809673653
$ws.Range("T9").Value = 0.221786809673653

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.055457
$ws.Range("H10").Value = 3.166371
$ws.Range("I10").Value = 0.009333819338610954
$ws.Range("J10").Value = 0.009333819338610952
$ws.Range("M10").Value = 3.795192333333334
$ws.Range("N10").Value = 11.385577
$ws.Range("O10").Value = 0.01044213755712683
$ws.Range("P10").Value = 0.01044213755712683
$ws.Range("Q10").Value = 4.005662314563
$ws.Range("R10").Value = 36.050960831067
$ws.Range("S10").Value = 0.00009746502546714618
$ws.Range("T10").Value = 0.00009746502546714616

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.055457
$ws.Range("H11").Value = 3.166371
$ws.Range("I11").Value = 0.009333819338610954
$ws.Range("J11").Value = 0.009333819338610952
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.6696287328350964
$ws.Range("P11").Value = 0.6696287328350964
$ws.Range("Q11").Value = 256.873323607524
$ws.Range("R11").Value = 2311.859912467716
$ws.Range("S11").Value = 0.00625019361622577
$ws.Range("T11").Value = 0.006250193616225769

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.055457
$ws.Range("H12").Value = 3.166371
$ws.Range("I12").Value = 0.009333819338610954
$ws.Range("J12").Value = 0.009333819338610952
$ws.Range("M12").Value = 29.801371
$ws.Range("N12").Value = 89.404113
$ws.Range("O12").Value = 0.08199584844219236
$ws.Range("P12").Value = 0.08199584844219235
$ws.Range("Q12").Value = 31.454065631547
$ws.Range("R12").Value = 283.086590683923
$ws.Range("S12").Value = 0.0007653344358755479
$ws.Range("T12").Value = 0.0007653344358755477

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.055457
$ws.Range("H13").Value = 3.166371
$ws.Range("I13").Value = 0.009333819338610954
$ws.Range("J13").Value = 0.009333819338610952
$ws.Range("M13").Value = 86.47679266666667
$ws.Range("N13").Value = 259.430378
$ws.Range("O13").Value = 0.2379332811655844
$ws.Range("P13").Value = 0.2379332811655844
$ws.Range("Q13").Value = 91.272536157582
$ws.Range("R13").Value = 821.452825418238
$ws.Range("S13").Value = 0.002220826261042489
$ws.Range("T13").Value = 0.002220826261042488

# Row 14
$ws.Range("G14").Value = 5.980707
$ws.Range("H14").Value = 17.942121
$ws.Range("I14").Value = 0.05288973274625675
$ws.Range("J14").Value = 0.05288973274625674
$ws.Range("M14").Value = 3.795192333333334
$ws.Range("N14").Value = 11.385577
$ws.Range("O14").Value = 0.01044213755712683
$ws.Range("P14").Value = 0.01044213755712683
$ws.Range("Q14").Value = 22.697933354313
$ws.Range("R14").Value = 204.281400188817
$ws.Range("S14").Value = 0.0005522818646960884
$ws.Range("T14").Value = 0.0005522818646960884

# Row 15
$ws.Range("G15").Value = 5.980707
$ws.Range("H15").Value = 17.942121
$ws.Range("I15").Value = 0.05288973274625675
$ws.Range("J15").Value = 0.05288973274625674
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.6696287328350964
$ws.Range("P15").Value = 0.6696287328350964
$ws.Range("Q15").Value = 1455.562931140524
$ws.Range("R15").Value = 13100.06638026472
$ws.Range("S15").Value = 0.03541648471886281
$ws.Range("T15").Value = 0.0354164847188628

# Row 16
$ws.Range("G16").Value = 5.980707
$ws.Range("H16").Value = 17.942121
$ws.Range("I16").Value = 0.05288973274625675
$ws.Range("J16").Value = 0.05288973274625674
$ws.Range("M16").Value = 29.801371
$ws.Range("N16").Value = 89.404113
$ws.Range("O16").Value = 0.08199584844219236
$ws.Range("P16").Value = 0.08199584844219235
$ws.Range("Q16").Value = 178.233268149297
$ws.Range("R16").Value = 1604.099413343673
$ws.Range("S16").Value = 0.004336738510410126
$ws.Range("T16").Value = 0.004336738510410125

# Row 17
$ws.Range("G17").Value = 5.980707
$ws.Range("H17").Value = 17.942121
$ws.Range("I17").Value = 0.05288973274625675
$ws.Range("J17").Value = 0.05288973274625674
$ws.Range("M17").Value = 86.47679266666667
$ws.Range("N17").Value = 259.430378
$ws.Range("O17").Value = 0.2379332811655844
$ws.Range("P17").Value = 0.2379332811655844
$ws.Range("Q17").Value = 517.192359239082
$ws.Range("R17").Value = 4654.731233151739
$ws.Range("S17").Value = 0.01258422765228772
$ws.Range("T17").Value = 0.01258422765228772

